$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Availability" column (G) with Digi-Key stock information ---

$header = "Availability as of 04-04-22"

$values = @(
  "In Stock",            # G2
  "In Stock",            # G3
  "In Stock",            # G4
  "In Stock",            # G5
  "In Stock",            # G6
  "Normally Stocking",   # G7
  "In Stock",            # G8
  "Normally Stocking",   # G9
  "In Stock",            # G10
  "In Stock",            # G11
  "In Stock",            # G12
  "In Stock",            # G13
  "In Stock",            # G14
  "In Stock",            # G15
  "In Stock",            # G16
  "In Stock",            # G17
  "In Stock",            # G18
  "Non-Stock",           # G19
  "Normally Stocking",   # G20
  "Normally Stocking",   # G21
  "Normally Stocking",   # G22
  "In Stock",            # G23
  "In Stock",            # G24
  "Non-Stock",           # G25
  "In Stock"             # G26
)

# Header cell - shares the same fill/border styling as the rest of row 1 (light grey fill,
# thin border) but without the wrap/vertical-center alignment used by A1:F1.
$g1 = $ws.Range("G1")
$g1.Value = $header
$g1.Interior.Color = 13882323

# Fill in all the data cells for the new column.
for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $cell = $ws.Cells.Item($row, 7)
  $cell.Value = $values[$i]
}

# Apply the thin black border to the whole new column (header + data) in one pass so it
# reuses the workbook's existing border style, then right-align just the data rows.
$full = $ws.Range("G1:G26")
$full.Borders.Color = 0
$full.Borders.LineStyle = 1

$dataRng = $ws.Range("G2:G26")
$dataRng.HorizontalAlignment = -4152

# Column width - fit to the new content
$ws.Columns("G").AutoFit()

# --- Conditional formatting: highlight any availability cell that is not "In Stock" ---

$rng = $ws.Range("G2:G26")
$fc = $rng.FormatConditions.Add(9, [System.Type]::Missing, [System.Type]::Missing, [System.Type]::Missing, "In Stock", 1)
$fc.Interior.Color = 5263615

# --- Update selection / scroll position to match the saved view ---

$ws.Range("K18").Select()
